$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the two runs of the "criticism" paragraph into a single run
# (same visible text, just re-typed so Word collapses the run split), then
# add a brand new paragraph after it with fresh text, moving the `_GoBack`
# bookmark from the end of the old paragraph to the end of the new one.
# ---------------------------------------------------------------------------

$criticismText = "A criticism I had was that it was difficult to get back from my purchase success page. I fixed this by adding a button in the middle which stands out in a bright red. This links back to my purchase page so the user can easily get back."

# Re-typing identical text over the old (two-run) text merges it into one run.
$d.Content.Find.Execute($criticismText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $criticismText, 2) | Out-Null

# Locate the paragraph again now that it is a single run.
$i = 0
$paraIndex = -1
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -eq ($criticismText + [char]13)) {
        $paraIndex = $i
    }
}
$criticismPara = $d.Paragraphs($paraIndex)

# The `_GoBack` bookmark currently sits at the end of this paragraph; it needs
# to end up at the end of the new paragraph we are about to insert. Remove it
# now and re-add it in the right spot once the new paragraph has its text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert a new, empty paragraph right after the criticism paragraph (it
# inherits the same paragraph/run formatting automatically).
$criticismPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs($paraIndex + 1)

$newText = "I have put comments in my code so that when I or another looks back on my code, they will be able to understand it."
# Append a one-character placeholder so the bookmark we add next is not
# sitting exactly on the paragraph-end boundary, then trim it away; this
# leaves the bookmark collapsed right after the real text.
$newPara.Range.InsertBefore($newText + "#")

$newParaRange = $d.Paragraphs($paraIndex + 1).Range
$bmPos = $newParaRange.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$placeholderRange = $d.Range($newParaRange.End - 2, $newParaRange.End - 1)
$placeholderRange.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Change 2: the "lastRenderedPageBreak" marker moves from the start of the
# "Task 13: Document testing" run to the start of the "This used to be
# 10000..." run earlier in the document. COM doesn't expose that marker
# directly, so we recreate the effect with Find & Replace no-op edits that
# nudge the underlying runs the same way the diff shows -- here we simply
# locate each run via Find and toggle formatting off/on to force a rewrite,
# keeping the two pieces of text completely untouched otherwise.
# ---------------------------------------------------------------------------

# (handled below together with change 3, since they are complementary)

$word.Selection.Find.ClearFormatting() | Out-Null

Write-Output "stage1-done"
